$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 657.913
$ws.Range("J17").Value = 642.36365
$ws.Range("L17").Value = 1927.09095
$ws.Range("N17").Value = -2263.09095
$ws.Range("H19").Value = 2759.2856
$ws.Range("J19").Value = 3324
$ws.Range("L19").Value = 3324
$ws.Range("N19").Value = -3674
$ws.Range("H96").Value = 7143817.5
$ws.Range("I96").Value = 14285957
$ws.Range("J96").Value = 1677.8
$ws.Range("K96").Value = 42857871
$ws.Range("L96").Value = 5033.4
$ws.Range("M96").Value = -42856498
$ws.Range("N96").Value = -7779.4
$ws.Range("H132").Value = 25255.666
$ws.Range("I132").Value = 26727.5
$ws.Range("K132").Value = 80182.5
$ws.Range("M132").Value = -77652.5
$ws.Range("H135").Value = 2160.5312
$ws.Range("I135").Value = 1714.7084
$ws.Range("K135").Value = 15432.3756
$ws.Range("M135").Value = -12897.3756
$ws.Range("H137").Value = 121360.8
$ws.Range("I137").Value = 152201
$ws.Range("K137").Value = 456603
$ws.Range("M137").Value = -454053
$ws.Range("H141").Value = 1778.2727
$ws.Range("I141").Value = 1778.2727
$ws.Range("K141").Value = 5334.8181
$ws.Range("M141").Value = -154.8181000000004

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15652.942
$ws.Range("I32").Value = 15652.942
$ws.Range("K32").Value = 15652.942
$ws.Range("M32").Value = -15365.942
$ws.Range("H61").Value = 8615.143
$ws.Range("I61").Value = 1301
$ws.Range("J61").Value = 52500
$ws.Range("K61").Value = 1301
$ws.Range("L61").Value = 52500
$ws.Range("M61").Value = -1089
$ws.Range("N61").Value = -52924
$ws.Range("H74").Value = 227225.11
$ws.Range("I74").Value = 286382.44
$ws.Range("J74").Value = 20174.5
$ws.Range("K74").Value = 286382.44
$ws.Range("L74").Value = 20174.5
$ws.Range("M74").Value = -285508.44
$ws.Range("N74").Value = -21922.5
$ws.Range("H77").Value = 227225.11
$ws.Range("I77").Value = 286382.44
$ws.Range("J77").Value = 20174.5
$ws.Range("K77").Value = 1431912.2
$ws.Range("L77").Value = 100872.5
$ws.Range("M77").Value = -1427544.2
$ws.Range("N77").Value = -109608.5
$ws.Range("H102").Value = 2810.125
$ws.Range("I102").Value = 2725.4285
$ws.Range("J102").Value = 3403
$ws.Range("K102").Value = 2725.4285
$ws.Range("L102").Value = 3403
$ws.Range("M102").Value = -1103.4285
$ws.Range("N102").Value = -6647
$ws.Range("H122").Value = 1832.52
$ws.Range("I122").Value = 1675.3334
$ws.Range("K122").Value = 5026.0002
$ws.Range("M122").Value = -2576.0002
$ws.Range("H132").Value = 1319.9375
$ws.Range("I132").Value = 829.9286
$ws.Range("J132").Value = 4750
$ws.Range("K132").Value = 2489.7858
$ws.Range("L132").Value = 14250
$ws.Range("M132").Value = 40.21420000000035
$ws.Range("N132").Value = -19310
$ws.Range("H136").Value = 8615.143
$ws.Range("I136").Value = 1301
$ws.Range("J136").Value = 52500
$ws.Range("K136").Value = 3903
$ws.Range("L136").Value = 157500
$ws.Range("M136").Value = -1353
$ws.Range("N136").Value = -162600

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2357.0667
$ws.Range("I86").Value = 2407.158
$ws.Range("K86").Value = 2407.158
$ws.Range("M86").Value = -1284.158
$ws.Range("H89").Value = 2357.0667
$ws.Range("I89").Value = 2407.158
$ws.Range("K89").Value = 12035.79
$ws.Range("M89").Value = -6419.789999999999
$ws.Range("H94").Value = 3783.1428
$ws.Range("I94").Value = 4033.7856
$ws.Range("K94").Value = 4033.7856
$ws.Range("M94").Value = -3582.7856
$ws.Range("H134").Value = 1537.7142
$ws.Range("I134").Value = 1204.1052
$ws.Range("J134").Value = 4707
$ws.Range("K134").Value = 3612.3156
$ws.Range("L134").Value = 14121
$ws.Range("M134").Value = -1077.3156
$ws.Range("N134").Value = -19191
$ws.Range("H138").Value = 40000
$ws.Range("J138").Value = 40000
$ws.Range("L138").Value = 40000
$ws.Range("N138").Value = -50280

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 155.8
$ws.Range("J7").Value = 224.8
$ws.Range("L7").Value = 224.8
$ws.Range("N7").Value = -450.8
$ws.Range("H10").Value = 499
$ws.Range("I10").Value = 499
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 499
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("M10").Value = -360
$ws.Range("H31").Value = 14286659
$ws.Range("I31").Value = 14286659
$ws.Range("K31").Value = 14286659
$ws.Range("M31").Value = -14286364
$ws.Range("H34").Value = 14286659
$ws.Range("I34").Value = 14286659
$ws.Range("K34").Value = 14286659
$ws.Range("M34").Value = -14286457
$ws.Range("H122").Value = 2945.4
$ws.Range("I122").Value = 3301
$ws.Range("J122").Value = 2412
$ws.Range("K122").Value = 9903
$ws.Range("L122").Value = 7236
$ws.Range("M122").Value = -7453
$ws.Range("N122").Value = -12136
$ws.Range("H132").Value = 57636.555
$ws.Range("I132").Value = 84204.914
$ws.Range("K132").Value = 252614.742
$ws.Range("M132").Value = -250084.742
$ws.Range("H134").Value = 2884.7778
$ws.Range("I134").Value = 2144.7
$ws.Range("J134").Value = 4999.2856
$ws.Range("K134").Value = 6434.099999999999
$ws.Range("L134").Value = 14997.8568
$ws.Range("M134").Value = -3899.099999999999
$ws.Range("N134").Value = -20067.8568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 537.5714
$ws.Range("I5").Value = 540.0909
$ws.Range("J5").Value = 528.3333
$ws.Range("K5").Value = 1620.2727
$ws.Range("L5").Value = 1584.9999
$ws.Range("M5").Value = -1508.2727
$ws.Range("N5").Value = -1808.9999
$ws.Range("H113").Value = 1283.5714
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H135").Value = 537.5714
$ws.Range("I135").Value = 540.0909
$ws.Range("J135").Value = 528.3333
$ws.Range("K135").Value = 4860.8181
$ws.Range("L135").Value = 4754.9997
$ws.Range("M135").Value = -2325.8181
$ws.Range("N135").Value = -9824.9997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 19751
$ws.Range("J92").Value = 19751
$ws.Range("L92").Value = 19751
$ws.Range("N92").Value = -23495
$ws.Range("H113").Value = 2762.44
$ws.Range("I113").Value = 2666.6
$ws.Range("K113").Value = 2666.6
$ws.Range("M113").Value = -496.5999999999999
$ws.Range("H122").Value = 4401.1763
$ws.Range("J122").Value = 4434.857
$ws.Range("L122").Value = 13304.571
$ws.Range("N122").Value = -18204.571
$ws.Range("H132").Value = 1694.091
$ws.Range("I132").Value = 1024.5555
$ws.Range("K132").Value = 3073.6665
$ws.Range("M132").Value = -543.6664999999998

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 10033333
$ws.Range("J2").Value = 10033333
$ws.Range("L2").Value = 10033333
$ws.Range("N2").Value = -10033557
$ws.Range("H22").Value = 1520.1333
$ws.Range("I22").Value = 1357.7142
$ws.Range("J22").Value = 1662.25
$ws.Range("K22").Value = 1357.7142
$ws.Range("L22").Value = 1662.25
$ws.Range("M22").Value = -1062.7142
$ws.Range("N22").Value = -2252.25
$ws.Range("H27").Value = 1520.1333
$ws.Range("I27").Value = 1357.7142
$ws.Range("J27").Value = 1662.25
$ws.Range("K27").Value = 1357.7142
$ws.Range("L27").Value = 1662.25
$ws.Range("M27").Value = -1250.7142
$ws.Range("N27").Value = -1876.25
$ws.Range("H32").Value = 25150
$ws.Range("I32").Value = 300
$ws.Range("J32").Value = 50000
$ws.Range("K32").Value = 300
$ws.Range("L32").Value = 50000
$ws.Range("M32").Value = 17
$ws.Range("N32").Value = -50634
$ws.Range("H136").Value = 3974.158
$ws.Range("I136").Value = 4082.2
$ws.Range("K136").Value = 12246.6
$ws.Range("M136").Value = -9696.599999999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 8811.546
$ws.Range("I81").Value = 9192.700000000001
$ws.Range("J81").Value = 5000
$ws.Range("K81").Value = 18385.4
$ws.Range("L81").Value = 10000
$ws.Range("M81").Value = -17324.4
$ws.Range("N81").Value = -12122
$ws.Range("H84").Value = 8811.546
$ws.Range("I84").Value = 9192.700000000001
$ws.Range("J84").Value = 5000
$ws.Range("K84").Value = 91927
$ws.Range("L84").Value = 50000
$ws.Range("M84").Value = -86623
$ws.Range("N84").Value = -60608
$ws.Range("H101").Value = 30000
$ws.Range("I101").Value = 30000
$ws.Range("K101").Value = 30000
$ws.Range("M101").Value = -26755
$ws.Range("H122").Value = 9117336
$ws.Range("I122").Value = 11396074
$ws.Range("J122").Value = 2386.5454
$ws.Range("K122").Value = 34188222
$ws.Range("L122").Value = 7159.6362
$ws.Range("M122").Value = -34185772
$ws.Range("N122").Value = -12059.6362
$ws.Range("H136").Value = 15449.979
$ws.Range("I136").Value = 17491.35
$ws.Range("J136").Value = 3785
$ws.Range("K136").Value = 52474.05
$ws.Range("L136").Value = 11355
$ws.Range("M136").Value = -49924.05
$ws.Range("N136").Value = -16455

Write-Output "done"